$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.480.58'
$ws.Range("E2").Value = '  -1.89%  '

$ws.Range("D3").Value = '2.896.73'
$ws.Range("E3").Value = '  -2.76%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.15'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.11'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.551'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.27%  '

$ws.Range("D9").Value = '2.906.64'
$ws.Range("E9").Value = '  -2.72%  '

$ws.Range("E10").Value = '  -5.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.01'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.359'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.62%  '

$ws.Range("D13").Value = '3.411.53'
$ws.Range("E13").Value = '  -2.52%  '

$ws.Range("E14").Value = '  +2.33%  '

$ws.Range("D15").Value = '60.507.75'
$ws.Range("E15").Value = '  -1.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.75'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.39%  '

$ws.Range("D17").Value = '2.909.04'
$ws.Range("E17").Value = '  -2.33%  '

$ws.Range("E18").Value = '  -4.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.98'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.65'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '358.50'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.65'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.75'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.453'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.26%  '

$ws.Range("E27").Value = '  -5.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.86'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.46%  '

$ws.Range("D30").Value = '0.0₃0840'
$ws.Range("E30").Value = '  -11.56%  '

$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("E32").Value = '  -2.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.72'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.91'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.33'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -8.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.58'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -7.52%  '

$ws.Range("E38").Value = '  -6.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.83'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.48'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.71'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.09%  '

$ws.Range("D42").Value = '2.290.60'
$ws.Range("E42").Value = '  -5.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.648'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0584'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.34'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.94'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0237'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.63%  '

$ws.Range("E49").Value = '  -1.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0917'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '248.83'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -8.23%  '
